# Auto-generated edit script to update '想去人数' (F column) values
# across the 4 worksheets per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3113
$ws.Range("F6").Value = 258
$ws.Range("F8").Value = 314
$ws.Range("F9").Value = 7309
$ws.Range("F10").Value = 65
$ws.Range("F12").Value = 54
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 423
$ws.Range("F17").Value = 1947
$ws.Range("F18").Value = 1777
$ws.Range("F21").Value = 67
$ws.Range("F22").Value = 1807
$ws.Range("F24").Value = 1226
$ws.Range("F25").Value = 637
$ws.Range("F27").Value = 1117
$ws.Range("F30").Value = 523
$ws.Range("F31").Value = 124
$ws.Range("F33").Value = 2680
$ws.Range("F34").Value = 1508
$ws.Range("F35").Value = 3004
$ws.Range("F36").Value = 2172
$ws.Range("F37").Value = 140
$ws.Range("F38").Value = 219
$ws.Range("F39").Value = 1149
$ws.Range("F41").Value = 41
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 372
$ws.Range("F46").Value = 243
$ws.Range("F48").Value = 739
$ws.Range("F49").Value = 433
$ws.Range("F50").Value = 106

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 537
$ws.Range("F20").Value = 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 2885
$ws.Range("F9").Value = 1125
$ws.Range("F10").Value = 1101
$ws.Range("F12").Value = 413
$ws.Range("F14").Value = 8124

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3113
$ws.Range("F5").Value = 258
$ws.Range("F7").Value = 2885
$ws.Range("F8").Value = 1125
$ws.Range("F9").Value = 1101
$ws.Range("F10").Value = 65
$ws.Range("F11").Value = 413
$ws.Range("F12").Value = 54
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 423
$ws.Range("F20").Value = 67
$ws.Range("F21").Value = 1807
$ws.Range("F23").Value = 1226
$ws.Range("F24").Value = 637
$ws.Range("F26").Value = 1117
$ws.Range("F30").Value = 537
$ws.Range("F31").Value = 523
$ws.Range("F32").Value = 124
$ws.Range("F34").Value = 2680
$ws.Range("F35").Value = 1508
$ws.Range("F36").Value = 3004
$ws.Range("F37").Value = 2172
$ws.Range("F38").Value = 140
$ws.Range("F39").Value = 219
$ws.Range("F40").Value = 1149
$ws.Range("F42").Value = 50
$ws.Range("F46").Value = 243
$ws.Range("F48").Value = 433
